# Updates cryptos list (Price / Volume(1h) columns) to the latest scraped
# values, and swaps the FraxShare/Aave row data (ranks 38/39 keep their
# position in column A, but the coin data that was in row 40 now belongs
# to row 41 and vice versa).
#
# Numeric-looking Price values are prefixed with a leading apostrophe so
# that Excel stores them as literal text (preserving formats such as
# "1.000" or "0.1220") instead of auto-converting them into numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.321.57'
$ws.Range("D3").Value = '1.932.55'
$ws.Range("E3").Value = '  +0.16%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = "'0.7559"
$ws.Range("E5").Value = '  +5.82%  '
$ws.Range("D6").Value = "'243.52"
$ws.Range("E6").Value = '  -2.14%  '
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").Value = "'27.99"
$ws.Range("E8").Value = '  +2.27%  '
$ws.Range("E9").Value = '  -0.99%  '
$ws.Range("D10").Value = "'0.07033"
$ws.Range("E10").Value = '  -0.89%  '
$ws.Range("D11").Value = "'0.7794"
$ws.Range("E11").Value = '  -1.59%  '
$ws.Range("E12").Value = '  -0.25%  '
$ws.Range("D13").Value = '1.942.56'
$ws.Range("E13").Value = '  +0.83%  '
$ws.Range("D14").Value = "'5.386"
$ws.Range("E14").Value = '  +0.21%  '
$ws.Range("D15").Value = "'93.27"
$ws.Range("E15").Value = '  -1.39%  '
$ws.Range("D16").Value = "'14.42"
$ws.Range("E16").Value = '  -1.74%  '
$ws.Range("D17").Value = '30.311.50'
$ws.Range("E17").Value = '  +0.13%  '
$ws.Range("D18").Value = "'252.68"
$ws.Range("D19").Value = "'5.982"
$ws.Range("E19").Value = '  +3.61%  '
$ws.Range("D20").Value = "'0.000007972"
$ws.Range("E20").Value = '  -0.68%  '
$ws.Range("D21").Value = '2.192.70'
$ws.Range("E21").Value = '  +0.52%  '
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("E24").Value = '  -1.58%  '
$ws.Range("D25").Value = "'9.512"
$ws.Range("E25").Value = '  -0.51%  '
$ws.Range("D26").Value = "'164.42"
$ws.Range("E26").Value = '  -0.98%  '
$ws.Range("D27").Value = "'19.06"
$ws.Range("E27").Value = '  -0.23%  '
$ws.Range("D28").Value = "'0.1313"
$ws.Range("E28").Value = '  +2.76%  '
$ws.Range("D29").Value = "'2.216"
$ws.Range("E29").Value = '  -2.72%  '
$ws.Range("D30").Value = "'1.370"
$ws.Range("E30").Value = '  +0.10%  '
$ws.Range("D31").Value = "'1.519"
$ws.Range("E31").Value = '  -0.62%  '
$ws.Range("D32").Value = "'4.405"
$ws.Range("E32").Value = '  +0.25%  '
$ws.Range("D33").Value = "'4.134"
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").Value = "'0.05224"
$ws.Range("E34").Value = '  +0.99%  '
$ws.Range("E35").Value = '  +4.13%  '
$ws.Range("D36").Value = "'0.7552"
$ws.Range("E36").Value = '  +1.30%  '
$ws.Range("D37").Value = "'2.794"
$ws.Range("E37").Value = '  +0.91%  '
$ws.Range("D38").Value = "'0.01947"
$ws.Range("E38").Value = '  -0.51%  '
$ws.Range("D39").Value = "'2.809"
$ws.Range("E39").Value = '  +0.05%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = "'6.616"
$ws.Range("E40").Value = '  +4.38%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").Value = "'78.67"
$ws.Range("E41").Value = '  +1.24%  '
$ws.Range("D42").Value = "'0.4486"
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("D43").Value = "'1.973"
$ws.Range("E43").Value = '  -0.37%  '
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = '  -0.14%  '
$ws.Range("D45").Value = "'0.8364"
$ws.Range("E45").Value = '  -1.11%  '
$ws.Range("D46").Value = "'9.931"
$ws.Range("E46").Value = '  +2.63%  '
$ws.Range("D47").Value = "'101.54"
$ws.Range("E47").Value = '  +0.83%  '
$ws.Range("D48").Value = "'7.576"
$ws.Range("E48").Value = '  +1.97%  '
$ws.Range("D49").Value = "'37.75"
$ws.Range("E49").Value = '  +3.70%  '
$ws.Range("D50").Value = "'983.69"
$ws.Range("E50").Value = '  +6.97%  '
$ws.Range("D51").Value = "'0.1220"
$ws.Range("E51").Value = '  +7.19%  '
